# Tag-1_6-Git-Workflow-im-Team.pptx
#   - "Brockhaus AG link von Folienmaster entfernt"
#   - "Formatfix in 2-3 Folie 6"
#
# Both edits live on the slide masters (the decorative "Rectangle 39..42"
# shapes and the cached date field both live on the master, not on an
# individual slide), so every Design/SlideMaster in the deck is inspected
# and the fix is applied wherever the relevant shapes actually exist.
#
# NOTE: shapes are looked up by scanning the Shapes collection by index and
# comparing .Name rather than calling Shapes.Item("Name") directly -- a
# failed name lookup on this host can otherwise resurface a stale shape
# reference left over from a previous lookup done through the same
# variable, so we avoid that pattern entirely.

$p = $ppt.ActivePresentation

for ($d = 1; $d -le $p.Designs.Count; $d++) {

    $master = $p.Designs.Item($d).SlideMaster

    $hlinkShape = $null
    $dateShape = $null

    for ($i = 1; $i -le $master.Shapes.Count; $i++) {
        $cand = $master.Shapes.Item($i)

        if ($cand.Name -eq "Rectangle 39") {
            $hlinkShape = $cand
        }
        if ($cand.Name -eq "Rectangle 6") {
            $dateShape = $cand
        }
    }

    # 1) Remove the rectangle carrying the external hyperlink to
    #    www.brockhaus-ag.de ("Rectangle 39"). The sibling helper
    #    rectangles (Rectangle 40/41/42) are left untouched and simply
    #    shift up one slot in the shape order.
    if ($hlinkShape -ne $null) {
        $hlinkShape.Delete()
    }

    # 2) Fix the cached date text from 13.06.2024 to 14.06.2024
    #    ("Rectangle 6" on the master holds the date field).
    if ($dateShape -ne $null) {
        if ($dateShape.HasTextFrame) {
            $tr = $dateShape.TextFrame.TextRange
            if ($tr.Text -eq "13.06.2024") {
                $tr.Text = "14.06.2024"
            }
        }
    }
}
